$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

function Get-ShapeByName($slide, $name) {
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $sh = $slide.Shapes.Item($i)
        if ($sh.Name -eq $name) { return $sh }
    }
    return $null
}

# Unicode punctuation used by the original deck.
$ldq = [char]0x201C   # “
$rdq = [char]0x201D   # ”
$lsq = [char]0x2018   # '
$rsq = [char]0x2019   # '

# --- Table 45 ------------------------------------------------------------
# Cell(2,1): "toAdd = “Name: David”" / "prevAddressBook = s2"
#         -> "toAdd = “Name: BTC”"   / "prevCoinBook = s2"
$shp = Get-ShapeByName $s "Table 45"
$tr = $shp.Table.Cell(2, 1).Shape.TextFrame.TextRange
$tr.Characters(23, 15).Text = "prevCoinBook"
$tr.Characters(9, 13).Text = $ldq + "Name: BTC" + $rdq

# --- Table 43 ------------------------------------------------------------
# Cell(2,1): "toAdd = “Name: David”" / "prevAddressBook = s2"
#         -> "toAdd = “Code: BTC”"   / "prevCoinBook = s2"
$shp = Get-ShapeByName $s "Table 43"
$tr = $shp.Table.Cell(2, 1).Shape.TextFrame.TextRange
$tr.Characters(23, 15).Text = "prevCoinBook"
$tr.Characters(9, 13).Text = $ldq + "Code: BTC" + $rdq

# --- Table 42 ------------------------------------------------------------
# Cell(2,1): "targetIndex = 5" / "prevAddressBook = s3" -> "prevCoinBook = s3"
$shp = Get-ShapeByName $s "Table 42"
$tr = $shp.Table.Cell(2, 1).Shape.TextFrame.TextRange
$tr.Characters(17, 15).Text = "prevCoinBook"

# --- Table 44 ------------------------------------------------------------
# Cell(2,1): "targetIndex = 5" / "prevAddressBook = s3" -> "prevCoinBook = s3"
$shp = Get-ShapeByName $s "Table 44"
$tr = $shp.Table.Cell(2, 1).Shape.TextFrame.TextRange
$tr.Characters(17, 15).Text = "prevCoinBook"

# --- TextBox 1 (caption under the diagram) --------------------------------
$shp = Get-ShapeByName $s "TextBox 1"
$tr = $shp.TextFrame.TextRange
$tr.Text = "The state of the coin book (before " + $lsq + "add c/BTC" + $rsq + " was executed) will be restored to s2."
